$wb = $excel.ActiveWorkbook

# Sheet "zh-cn": Correspond Handoff Datetime (E2) and Correspond Handback DateTime (H2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-21 02:55:35"
$wsZhCn.Range("H2").Value = "2016-03-21 02:55:54"

# Sheet "de-de": Correspond Handoff Datetime (E2) and Correspond Handback DateTime (H2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-21 02:55:38"
$wsDeDe.Range("H2").Value = "2016-03-21 02:56:00"
